$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2 and G3 values (both relays re-set to 600 A)
$ws.Range("G2").Value = 600
$ws.Range("G3").Value = 600

# The case for row 4 (relay #3) was "executed" - its entered parameters
# are removed again. Cells that hold a non-default (per-column) style
# (A4, C4:F4, H4) keep their formatting but lose their value, matching
# the still-blank template rows below (5-9). Cells that only carried the
# plain column-default style (B4, G4, I4:V4) are cleared completely, so
# they disappear from the sheet just like their counterparts in rows 5-9.
$ws.Range("A4").ClearContents()
$ws.Range("C4:F4").ClearContents()
$ws.Range("H4").ClearContents()

$ws.Range("B4").Clear()
$ws.Range("G4").Clear()
$ws.Range("I4:V4").Clear()

# Move the active selection to G4
$ws.Range("G4").Select()
